# Apply updated cryptocurrency price / volume figures to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    # Force the cell to remain a text value (matches the source data, which stores
    # these look-like-numbers strings as text) and then drop back to the default
    # "Normal" style so no stray number-format is left attached to the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.719.64"
$ws.Range("E2").Value = "  +0.58%  "
Set-TextValue $ws.Range("D3") "1.639.39"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextValue $ws.Range("D5") "212.72"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  -2.15%  "
$ws.Range("E7").Value = "  -0.06%  "
Set-TextValue $ws.Range("D8") "23.25"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("E10").Value = "  +0.02%  "
Set-TextValue $ws.Range("D11") "0.0890"
$ws.Range("E11").Value = "  +0.07%  "
Set-TextValue $ws.Range("D12") "1.872.00"
$ws.Range("E12").Value = "  -0.45%  "
Set-TextValue $ws.Range("D13") "1.645.09"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("E14").Value = "  +0.52%  "
Set-TextValue $ws.Range("D15") "0.562"
$ws.Range("E15").Value = "  -3.72%  "
Set-TextValue $ws.Range("D16") "64.81"
$ws.Range("E16").Value = "  +0.50%  "
Set-TextValue $ws.Range("D17") "27.687.02"
$ws.Range("E17").Value = "  +0.58%  "
Set-TextValue $ws.Range("D18") "230.38"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  +2.25%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +0.13%  "
Set-TextValue $ws.Range("D23") "10.27"
$ws.Range("E23").Value = "  +5.55%  "
Set-TextValue $ws.Range("D24") "2.08"
$ws.Range("E24").Value = "  +3.21%  "
Set-TextValue $ws.Range("D25") "150.97"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("E28").Value = "  -0.06%  "
Set-TextValue $ws.Range("D29") "15.60"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("E30").Value = "  +0.40%  "
Set-TextValue $ws.Range("D31") "0.0487"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +0.24%  "
Set-TextValue $ws.Range("D33") "1.461.33"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("E36").Value = "  -0.41%  "
Set-TextValue $ws.Range("D37") "0.568"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("E39").Value = "  +0.40%  "
Set-TextValue $ws.Range("D40") "0.895"
$ws.Range("E40").Value = "  +9.63%  "
Set-TextValue $ws.Range("D41") "69.08"
$ws.Range("E41").Value = "  +6.14%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("E44").Value = "  +1.21%  "
$ws.Range("E45").Value = "  -0.72%  "
Set-TextValue $ws.Range("D47") "1.781.94"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("E48").Value = "  +3.41%  "
Set-TextValue $ws.Range("D49") "87.06"
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("E50").Value = "  -1.18%  "
